# Atualizei dados para BIBI e ADD 06-05-2025
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12: cohort 2022, period_index 4 -> num_customers 66 -> 67, retention_rate recalculated (67/419)
$ws.Range("C12").Value = 67
$ws.Range("E12").Value = 0.1599045346062052

# Row 19: cohort 2023, period_index 2 -> num_customers 61 -> 64, retention_rate recalculated (64/123)
$ws.Range("C19").Value = 64
$ws.Range("E19").Value = 0.5203252032520326

# Row 21: cohort 2024, period_index 1 -> num_customers 107 -> 108, retention_rate recalculated (108/206)
$ws.Range("C21").Value = 108
$ws.Range("E21").Value = 0.5242718446601942

# Row 22: cohort 2025, period_index 0 -> num_customers 35 -> 42, cohort_size 35 -> 42, retention_rate stays 1
$ws.Range("C22").Value = 42
$ws.Range("D22").Value = 42
